$d = $word.ActiveDocument

# Turn off track-changes so our edits land as plain content, not revisions.
$d.TrackRevisions = $false

# ---------------------------------------------------------------------
# Step 1: normalize run formatting across the whole (still single)
# paragraph so every run has the complex-script font/size set
# (w:cs="SimSun", w:szCs="24") to match the target markup.
# ---------------------------------------------------------------------
$whole = $d.Paragraphs(1).Range
$whole.Font.NameBi = "SimSun"
$whole.Font.SizeBi = 12

# ---------------------------------------------------------------------
# Step 2: remove the trailing space before the final period-less run
# "...管理硕士学位。 " -> "...管理硕士学位。"
# ---------------------------------------------------------------------
$trailingSpace = $d.Content
$trailingSpace.Find.ClearFormatting()
$found = $trailingSpace.Find.Execute("学位。 ", $false, $false, $false, $false, $false, $true, 1, $false, "学位。", 2)

# ---------------------------------------------------------------------
# Step 3: split the single paragraph into three content paragraphs by
# inserting paragraph breaks at the two internal boundaries, and add a
# blank paragraph after each break (giving 5 paragraphs total).
# ---------------------------------------------------------------------

# Boundary A: after "...Equifax 在全球 14 个国家的业务。" and before "Alan 曾担任..."
$rngA = $d.Content
$rngA.Find.ClearFormatting()
$foundA = $rngA.Find.Execute("的业务。", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngA.Collapse(0)
$rngA.InsertParagraphAfter()
$rngA.Collapse(0)
$rngA.MoveStart(1, 1) | Out-Null
$rngA.InsertParagraphAfter()

# Boundary B: after "...贡献巨大。" and before "除了获得多个资格证书..."
$rngB = $d.Content
$rngB.Find.ClearFormatting()
$foundB = $rngB.Find.Execute("贡献巨大。", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngB.Collapse(0)
$rngB.InsertParagraphAfter()
$rngB.Collapse(0)
$rngB.MoveStart(1, 1) | Out-Null
$rngB.InsertParagraphAfter()

# ---------------------------------------------------------------------
# Step 4: apply first-line indent (720 twips = 36pt) to the three
# content paragraphs (1st, 3rd and 5th); the two blank spacer
# paragraphs (2nd and 4th) keep no indent.
# ---------------------------------------------------------------------
$d.Paragraphs(1).Format.FirstLineIndent = 36
$d.Paragraphs(3).Format.FirstLineIndent = 36
$d.Paragraphs(5).Format.FirstLineIndent = 36

# ---------------------------------------------------------------------
# Step 5: add the _GoBack bookmark at the start of the final paragraph.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs(5)
$bookmarkRange = $lastPara.Range.Duplicate
$bookmarkRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
